$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and two swapped rows)
# D-column numeric-looking values are forced to remain plain text (matching the
# original inlineStr string cells) by temporarily marking the cell as Text format,
# then restoring the number format back to General/Normal so no lingering style
# attribute is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.584.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.923.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.75%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("E6").Value = '  +0.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4820'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4052'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08186'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.008'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.898.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.055'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.276'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06872'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.37%  '
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001036'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.011'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.589.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.637'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.201'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.150.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.32'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.345'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.71%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.99'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.081'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.47'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.001'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09591'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.594'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.389'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06501'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02281'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.211'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5917'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.14%  '
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.845'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.514'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1840'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.282'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07529'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.29%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5535'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.953'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '117.57'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.425'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.54%  '
